$d = $word.ActiveDocument
$needle = "пометка ячейки как пустая"
$rng = $d.Content.Duplicate
$ok = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
Write-Host "ok=$ok"
